$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "You are an AI product manager for an e-commerce platform and want to enhance the recommendation engine using AI. Your goal is to improve product recommendations for users based on their past behavior and preferences. What should be your initial step in this AI-driven recommendation project?",
        "ques_type": 2,
        "options": [
            "Analyze historical user interaction data.",
            "Collect user feedback on current recommendations.",
            "Train a neural network for recommendation modeling.",
            "Implement a new user interface design."
        ],
        "score": "Analyze historical user interaction data."
    },
    {
        "title": "You are an AI product manager and have identified a promising AI use case for your company's new chatbot product: personalized product recommendations based on user preferences. What should be your immediate action to develop this AI feature successfully?",
        "ques_type": 2,
        "options": [
            "Conduct a thorough analysis of user data and preferences to inform the AI recommendation model.",
            "Begin coding the AI recommendation algorithm to see how it performs in a live environment.",
            "Design an appealing user interface for the chatbot's recommendations without assessing user data.",
            "Collaborate with marketing to create a promotional campaign for the AI recommendation feature."
        ],
        "score": "Conduct a thorough analysis of user data and preferences to inform the AI recommendation model."
    },
    {
        "title": "You are an AI product manager working on a healthcare AI system that uses patient data for disease prediction. You've just discovered a bias in the training data that could lead to inaccurate predictions for certain demographic groups. What should you do to address this issue?",
        "ques_type": 2,
        "options": [
            "Develop a bias mitigation strategy and retrain the model with fairer data.",
            "Retrain the model with the biased data because it may improve overall accuracy.",
            "Ignore the bias since it's challenging to rectify biases in AI systems.",
            "Continue using the current model bias is an inherent part of AI."
        ],
        "score": "Develop a bias mitigation strategy and retrain the model with fairer data."
    },
    {
        "title": "You are the AI product manager for a healthcare company working on an AI-driven diagnostic tool. Your development team comprises data scientists, software engineers, and healthcare domain experts. You've received feedback that the AI tool's predictions are not consistently accurate. What is your immediate course of action to address this issue?",
        "ques_type": 2,
        "options": [
            "Review of the model training data and evaluation metrics to identify potential bias or data-quality issues.",
            "Instruct the software engineers to build a new AI model to produce more accurate predictions.",
            "Initiate a complete redesign of the AI tool\u2019s user interface and user experience.",
            "Replace the data science team with a new group of experts in AI."
        ],
        "score": "Review of the model training data and evaluation metrics to identify potential bias or data-quality issues."
    }
]
'@

# The single JSON/"questions" blob used to live in A2 (A1 held a throwaway
# "0" placeholder with bold+boxed formatting). Drop A1's old value/format,
# pull the reformatted text up into A1, and delete the now-empty row 2 so
# the sheet is back down to a single A1 cell.
$ws.Range("A1").ClearFormats() | Out-Null
$ws.Range("A1").Value = $newText
$ws.Range("A2").EntireRow.Delete() | Out-Null

# Setting a long, multi-line value auto-expands the row to a custom height;
# AutoFit() recomputes it back to a non-custom (default-tracking) height so
# no stray ht=/customHeight= survives in the saved XML.
$ws.Rows.Item(1).AutoFit() | Out-Null
